$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same
#    column layout / header wording / styles), placed right before the
#    "总计" summary sheet.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Row 2: 516980 / 华富中证证券公司先锋策略ETF / 0.39 / 98.42 / 2.49 / 0.0097 / 10
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "516980"
$newSheet.Range("C2").Value = "华富中证证券公司先锋策略ETF"
$newSheet.Range("D2").Value = "0.39"
$newSheet.Range("E2").Value = "98.42"
$newSheet.Range("F2").Value = "2.49"
$newSheet.Range("G2").Value = "0.0097"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 10

# Row 3 (new): 090011 / 大成核心双动力混合 / 0.34 / 93.14 / 2.07 / 0.0070 / 6
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A3:H3").PasteSpecial(-4122)
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "090011"
$newSheet.Range("C3").Value = "大成核心双动力混合"
$newSheet.Range("D3").Value = "0.34"
$newSheet.Range("E3").Value = "93.14"
$newSheet.Range("F3").Value = "2.07"
$newSheet.Range("G3").Value = "0.0070"
$newSheet.Range("B3:G3").ClearFormats()
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q1,
#    pushing the existing quarters down by one.
#    (Re-fetch by name: the old $totalSheet handle now resolves to the
#    sheet that took its former tab position, i.e. the copy we just made.)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.02

# The leading "A" column is a running 0-based index; renumber the rows
# that got pushed down by the insert (were 0,1,2 -> now 1,2,3).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the original active sheet/selection.
$wb.Worksheets.Item("2021-Q1").Activate()
